$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.443.65"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.737.31"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4546"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +7.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3526"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07378"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.075"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.39"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.911"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.046"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").Value = "1.732.74"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.05"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001053"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06331"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.750"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").Value = "27.481.58"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.073"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.46"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "1.928.87"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.047"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.90"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.044"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09098"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.649"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.392"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02265"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.59"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05951"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2052"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6222"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.878"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.186"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.370"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5785"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.83"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.923"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06840"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.111"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.94%  "
